# Regenerate the "K" (strikeouts) column (column G) of save_data using
# strikeout counts (K) instead of the previous "Strike#" pitch-count based
# values, as part of regenerating std/mean and s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new K value (column G), keyed by the sheet row number (data starts
# at row 2; row 1 is the header row).
$kValues = @{
    2  = 4
    3  = 9
    4  = 5
    5  = 12
    6  = 5
    7  = 5
    8  = 5
    9  = 2
    10 = 7
    11 = 8
    12 = 2
    13 = 1
    14 = 4
    15 = 3
    16 = 10
    17 = 5
    18 = 4
    19 = 6
    20 = 8
    21 = 5
    22 = 6
    23 = 1
    24 = 1
    25 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
